$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2
$ws.Range("C2").Value = 5211018942
$ws.Range("D2").Value = "Карпеев Данил Алексеевич"
$ws.Range("E2").Value = "Параша"

# Update existing row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "улгту"
$ws.Range("D3").Value = "привет"
$ws.Range("E3").Value = "привет"

# Add new row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "привет"
$ws.Range("C4").Value = 5211018942
$ws.Range("D4").Value = "Карпеев Данил Алексеевич"
$ws.Range("E4").Value = "Хорошая"

# Add new row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "привет"
$ws.Range("C5").Value = 5813154625
$ws.Range("D5").Value = "привет"
$ws.Range("E5").Value = "салам всем мусорам"
